# Refresh cryptos list values (prices / 1h volume %) and reorder a few rows
# per the latest coinranking.com snapshot pulled by the scheduled GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.251.31'
$ws.Range("E2").Value = '  -0.76%  '

$ws.Range("D3").Value = '1.880.01'
$ws.Range("E3").Value = '  -1.74%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''234.89'
$ws.Range("E5").Value = '  -1.73%  '

$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").Value = '''0.4678'
$ws.Range("E7").Value = '  -2.12%  '

$ws.Range("D8").Value = '''0.2824'
$ws.Range("E8").Value = '  -0.66%  '

$ws.Range("D9").Value = '''0.06585'
$ws.Range("E9").Value = '  -1.61%  '

$ws.Range("D10").Value = '''20.58'
$ws.Range("E10").Value = '  +9.37%  '

$ws.Range("D11").Value = '''0.07764'
$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").Value = '''97.41'
$ws.Range("E12").Value = '  -4.47%  '

$ws.Range("D13").Value = '1.879.57'
$ws.Range("E13").Value = '  -1.70%  '

$ws.Range("D14").Value = '''5.064'
$ws.Range("E14").Value = '  -2.83%  '

$ws.Range("D15").Value = '''0.6722'
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("D16").Value = '''283.49'
$ws.Range("E16").Value = '  +5.62%  '

$ws.Range("D17").Value = '30.255.35'
$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").Value = '''0.9997'
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").Value = '''12.59'
$ws.Range("E19").Value = '  -0.60%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.125.03'
$ws.Range("E20").Value = '  -1.45%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''5.388'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '''0.000007244'
$ws.Range("E22").Value = '  -2.92%  '

$ws.Range("D23").Value = '''1.002'

$ws.Range("D24").Value = '''6.159'
$ws.Range("E24").Value = '  -2.13%  '

$ws.Range("D25").Value = '''9.345'
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("D26").Value = '''167.78'
$ws.Range("E26").Value = '  +0.39%  '

$ws.Range("D27").Value = '''19.14'
$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("D28").Value = '''1.982'
$ws.Range("E28").Value = '  -3.84%  '

$ws.Range("D29").Value = '''1.382'
$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("D30").Value = '''0.09662'
$ws.Range("E30").Value = '  -3.29%  '

$ws.Range("D31").Value = '''4.362'
$ws.Range("E31").Value = '  -6.78%  '

$ws.Range("D32").Value = '''1.472'
$ws.Range("E32").Value = '  -2.49%  '

$ws.Range("D33").Value = '''4.104'
$ws.Range("E33").Value = '  -3.28%  '

$ws.Range("D34").Value = '''0.04660'
$ws.Range("E34").Value = '  -1.33%  '

$ws.Range("D35").Value = '''0.7037'
$ws.Range("E35").Value = '  -3.17%  '

$ws.Range("D36").Value = '''1.093'
$ws.Range("E36").Value = '  -1.27%  '

$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '''1.001'
$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '''2.714'
$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01866'
$ws.Range("E39").Value = '  -2.51%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''6.580'
$ws.Range("E40").Value = '  +5.65%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.523'
$ws.Range("E41").Value = '  -3.73%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''71.83'
$ws.Range("E42").Value = '  -4.12%  '

$ws.Range("D43").Value = '''0.8625'
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''1.954'
$ws.Range("E44").Value = '  -0.74%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''1.001'
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '''102.89'
$ws.Range("E46").Value = '  -2.85%  '

$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").Value = '''0.4172'
$ws.Range("E47").Value = '  -1.93%  '

$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '''981.64'
$ws.Range("E48").Value = '  +7.39%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '''7.242'
$ws.Range("E49").Value = '  -2.17%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.147'
$ws.Range("E50").Value = '  +4.57%  '

$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '''33.81'
$ws.Range("E51").Value = '  -2.55%  '
